$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in column H, matching the formatting of the other
# header cells (bold, centered, bordered - same style as G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data value for the "Save" column on row 2.
$ws.Range("H2").Value = 0
